$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Update table cell text: "Tipo de Sequenciamento" -> "Tipo de amostra"
$tbl = $s.Shapes.Item(2).Table
$cell = $tbl.Cell(1, 4)
$cell.Shape.TextFrame.TextRange.Runs(2).Text = "amostra"

# Nudge the picture (id=27, "Picture 26") to its new position.
$pic = $s.Shapes.Item(3)
$pic.Left = 9323846 / 12700
$pic.Top = 5041771 / 12700
